$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9 corresponds to @nit_srinagar_2023_batch
$ws.Range("E9").Value = "2026-02-13T04:57:49.493927+00:00"
$ws.Range("H9").Value = 1
$ws.Range("I9").Value = 2
$ws.Range("L9").Value = "[33]"
$ws.Range("M9").Value = "[24, 19]"
